$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "0.9991").
# Force them to remain plain text (matching the source inlineStr cells)
# by temporarily applying a text number format, then restore the default
# "Normal" style so no stray style index is left on the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.090.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -2.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.849.59'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.19%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9991'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.6940'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -5.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '237.79'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.39%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9985'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07777'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +9.19%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3044'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.92%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.31'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.47%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08112'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.56%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.857.11'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.7258'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.17%  '

$ws.Range("E14").Value = '  -2.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '88.99'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.67%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.096.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.46%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.744'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.37%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007821'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.19'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.56%  '

$ws.Range("E20").Value = '  -4.92%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9992'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.28%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.098.93'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.83%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9984'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.25%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.598'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.00%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.978'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.42%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '161.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.84%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1432'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.03%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.980'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.95%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.400'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.32%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.487'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.70%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.486'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.30%  '

$ws.Range("E33").Value = '  -4.15%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05229'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.45%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.182'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.01%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7042'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.02%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.003'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.643'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.14%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01855'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.13%  '

$ws.Range("E40").Value = '  -2.45%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9146'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.39%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.094.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.71%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.990'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.34%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.4266'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.40%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '70.66'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9978'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.26%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '102.85'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.04%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.769'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.71%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.995.44'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.55%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.171'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.58%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.981'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.02%  '
